$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 72.79331566666666
$ws.Range("H2").Value = 218.379947
$ws.Range("I2").Value = 0.2828741606141505
$ws.Range("J2").Value = 0.2828741606141506
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.8749903333333333
$ws.Range("N2").Value = 2.624971
$ws.Range("O2").Value = 0.2670516933349977
$ws.Range("P2").Value = 0.2670516933349977
$ws.Range("Q2").Value = 63.69344753961521
$ws.Range("R2").Value = 573.2410278565369
$ws.Range("S2").Value = 0.07554202359272502
$ws.Range("T2").Value = 0.07554202359272504

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 72.79331566666666
$ws.Range("H3").Value = 218.379947
$ws.Range("I3").Value = 0.2828741606141505
$ws.Range("J3").Value = 0.2828741606141506
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.8147036666666666
$ws.Range("N3").Value = 2.444111
$ws.Range("O3").Value = 0.2486518827250642
$ws.Range("P3").Value = 0.2486518827250642
$ws.Range("Q3").Value = 59.30498118245743
$ws.Range("R3").Value = 533.7448306421169
$ws.Range("S3").Value = 0.07033719261098074
$ws.Range("T3").Value = 0.07033719261098076

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 72.79331566666666
$ws.Range("H4").Value = 218.379947
$ws.Range("I4").Value = 0.2828741606141505
$ws.Range("J4").Value = 0.2828741606141506
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.586789
$ws.Range("N4").Value = 4.760367
$ws.Range("O4").Value = 0.484296423939938
$ws.Range("P4").Value = 0.484296423939938
$ws.Range("Q4").Value = 115.5076325733943
$ws.Range("R4").Value = 1039.568693160549
$ws.Range("S4").Value = 0.1369949444104448
$ws.Range("T4").Value = 0.1369949444104448

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 117.1700846666667
$ws.Range("H5").Value = 351.510254
$ws.Range("I5").Value = 0.4553218801152877
$ws.Range("J5").Value = 0.4553218801152878
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.8749903333333333
$ws.Range("N5").Value = 2.624971
$ws.Range("O5").Value = 0.2670516933349977
$ws.Range("P5").Value = 0.2670516933349977
$ws.Range("Q5").Value = 102.5226914391816
$ws.Range("R5").Value = 922.704222952634
$ws.Range("S5").Value = 0.1215944790972624
$ws.Range("T5").Value = 0.1215944790972624

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 117.1700846666667
$ws.Range("H6").Value = 351.510254
$ws.Range("I6").Value = 0.4553218801152877
$ws.Range("J6").Value = 0.4553218801152878
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.8147036666666666
$ws.Range("N6").Value = 2.444111
$ws.Range("O6").Value = 0.2486518827250642
$ws.Range("P6").Value = 0.2486518827250642
$ws.Range("Q6").Value = 95.45889760157712
$ws.Range("R6").Value = 859.130078414194
$ws.Range("S6").Value = 0.1132166427365823
$ws.Range("T6").Value = 0.1132166427365823

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 117.1700846666667
$ws.Range("H7").Value = 351.510254
$ws.Range("I7").Value = 0.4553218801152877
$ws.Range("J7").Value = 0.4553218801152878
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.586789
$ws.Range("N7").Value = 4.760367
$ws.Range("O7").Value = 0.484296423939938
$ws.Range("P7").Value = 0.484296423939938
$ws.Range("Q7").Value = 185.9242014781353
$ws.Range("R7").Value = 1673.317813303218
$ws.Range("S7").Value = 0.220510758281443
$ws.Range("T7").Value = 0.220510758281443

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 67.37122333333333
$ws.Range("H8").Value = 202.11367
$ws.Range("I8").Value = 0.2618039592705617
$ws.Range("J8").Value = 0.2618039592705618
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.8749903333333333
$ws.Range("N8").Value = 2.624971
$ws.Range("O8").Value = 0.2670516933349977
$ws.Range("P8").Value = 0.2670516933349977
$ws.Range("Q8").Value = 58.94916916150778
$ws.Range("R8").Value = 530.54252245357
$ws.Range("S8").Value = 0.06991519064501028
$ws.Range("T8").Value = 0.0699151906450103

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 67.37122333333333
$ws.Range("H9").Value = 202.11367
$ws.Range("I9").Value = 0.2618039592705617
$ws.Range("J9").Value = 0.2618039592705618
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.8147036666666666
$ws.Range("N9").Value = 2.444111
$ws.Range("O9").Value = 0.2486518827250642
$ws.Range("P9").Value = 0.2486518827250642
$ws.Range("Q9").Value = 54.88758267748555
$ws.Range("R9").Value = 493.98824409737
$ws.Range("S9").Value = 0.0650980473775012
$ws.Range("T9").Value = 0.06509804737750122

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 67.37122333333333
$ws.Range("H10").Value = 202.11367
$ws.Range("I10").Value = 0.2618039592705617
$ws.Range("J10").Value = 0.2618039592705618
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.586789
$ws.Range("N10").Value = 4.760367
$ws.Range("O10").Value = 0.484296423939938
$ws.Range("P10").Value = 0.484296423939938
$ws.Range("Q10").Value = 106.9039161018766
$ws.Range("R10").Value = 962.13524491689
$ws.Range("S10").Value = 0.1267907212480502
$ws.Range("T10").Value = 0.1267907212480502
